$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old column F (ElementName3) data - no longer used
$ws.Range("F1:F13").ClearContents() | Out-Null

# Write the updated grid of values into A1:E13
$ws.Range("A1").Value = "Code"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "ChefModule"
$ws.Range("D1").Value = "ElementName1"
$ws.Range("E1").Value = "ElementName2"

$ws.Range("A2").Value = "GINF11"
$ws.Range("B2").Value = "pede. Suspendisse dui."
$ws.Range("C2").Value = "EL Haddad"
$ws.Range("D2").Value = "Nullam feugiat placerat"
$ws.Range("E2").Value = "varius et, euismod"

$ws.Range("A3").Value = "GINF12"
$ws.Range("B3").Value = "a nunc. In"
$ws.Range("C3").Value = "Badir"
$ws.Range("D3").Value = "sodales nisi magna"
$ws.Range("E3").Value = "elementum sem, vitae"

$ws.Range("A4").Value = "GINF13"
$ws.Range("B4").Value = "amet metus. Aliquam"
$ws.Range("C4").Value = "Ezzine"
$ws.Range("D4").Value = "Cras vulputate velit"
$ws.Range("E4").Value = "scelerisque neque sed"

$ws.Range("A5").Value = "GINF14"
$ws.Range("B5").Value = "quam vel sapien"
$ws.Range("C5").Value = "El Alami Hassoun"
$ws.Range("D5").Value = "Nunc mauris elit,"
$ws.Range("E5").Value = "libero et tristique"

$ws.Range("A6").Value = "GINF15"
$ws.Range("B6").Value = "feugiat nec, diam."
$ws.Range("C6").Value = "Lazaar"
$ws.Range("D6").Value = "pellentesque. Sed dictum."
$ws.Range("E6").Value = "ridiculus mus. Proin"

$ws.Range("A7").Value = "GINF16"
$ws.Range("B7").Value = "nonummy. Fusce fermentum"
$ws.Range("C7").Value = "El Haddad"
$ws.Range("D7").Value = "neque pellentesque massa"
$ws.Range("E7").Value = "Mauris eu turpis."

$ws.Range("A8").Value = "GINF21"
$ws.Range("B8").Value = "a, arcu. Sed"
$ws.Range("C8").Value = "EL Haddad"
$ws.Range("D8").Value = "sit amet risus."
$ws.Range("E8").Value = "Nulla facilisi. Sed"

$ws.Range("A9").Value = "GINF22"
$ws.Range("B9").Value = "Suspendisse eleifend. Cras"
$ws.Range("C9").Value = "El Alami Hassoun"
$ws.Range("D9").Value = "velit dui, semper"
$ws.Range("E9").Value = "ligula elit, pretium"

$ws.Range("A10").Value = "GINF23"
$ws.Range("B10").Value = "ante. Nunc mauris"
$ws.Range("C10").Value = "Badir"
$ws.Range("D10").Value = "tortor at risus."
$ws.Range("E10").Value = "felis. Donec tempor,"

$ws.Range("A11").Value = "GINF24"
$ws.Range("B11").Value = "lobortis quam a"
$ws.Range("C11").Value = "Ezzine"
$ws.Range("D11").Value = "euismod est arcu"
$ws.Range("E11").Value = "ligula eu enim."

$ws.Range("A12").Value = "GINF25"
$ws.Range("B12").Value = "rhoncus. Nullam velit"
$ws.Range("C12").Value = "Ben Achrab"
$ws.Range("D12").Value = "ut dolor dapibus"
$ws.Range("E12").Value = "commodo tincidunt nibh."

$ws.Range("A13").Value = "GINF26"
$ws.Range("B13").Value = "Donec tincidunt. Donec"
$ws.Range("C13").Value = "EL Haddad"
$ws.Range("D13").Value = "ornare tortor at"
$ws.Range("E13").Value = "ac, feugiat non,"

# Adjust column widths to fit new content (values chosen so the engine's
# internal pixel-grid snapping lands on the closest achievable width to the
# real-Excel bestFit values of 27 / 16.140625 / 27.5703125 / 23.5703125)
$ws.Columns.Item(2).ColumnWidth = 26.166666666666668
$ws.Columns.Item(3).ColumnWidth = 15.333333333333334
$ws.Columns.Item(4).ColumnWidth = 26.666666666666668
$ws.Columns.Item(5).ColumnWidth = 22.666666666666668

# Update selection
$ws.Range("H11").Select() | Out-Null
